$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before G. This shifts the old column G (header
# "2035_Fx" plus its data) over to column H, while the freshly inserted
# (blank) column G inherits number formatting from its left neighbour F.
$ws.Columns("G").Insert()

# Fill the two pre-existing blank placeholder columns (C, D) with the new
# "2020" and "2025nb" scenario data. D2/D3 lose their number formatting in
# the source workbook, so clear formats on just those two cells first.
$ws.Range("C1").Value = 2020
$ws.Range("C2").Value = 82
$ws.Range("C3").Value = 7.8562630000000002
$ws.Range("C4").Value = 7.0709140000000001

$ws.Range("D1").Value = "2025nb"
$ws.Range("D2").ClearFormats()
$ws.Range("D2").Value = 81
$ws.Range("D3").ClearFormats()
$ws.Range("D3").Value = 7.7965619999999998
$ws.Range("D4").Value = 7.0202210000000003

# Update the shifted former-G (now H) column header/value to its new text.
$ws.Range("H1").Value = "2035_F                                        "

# Fill the newly inserted "2035 E" column.
$ws.Range("G1").Value = "2035 E"
$ws.Range("G2").Value = 83
$ws.Range("G3").Value = 7.3506309999999999
$ws.Range("G4").Value = 6.8323330000000002

# Two new blank formatted rows below the existing placeholder rows.
$ws.Range("A15:B16").NumberFormat = "0.00"

$ws.Range("G4").Select()
